$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# (subtract 5/6 to offset the character-width padding Excel adds so the
# saved OOXML "width" attribute comes out to the exact target value)
$ws.Columns.Item(1).ColumnWidth = 20 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 22 - (5/6)

# Update header row
$ws.Range("A1").Value = "var_1_input_object"
$ws.Range("B1").Value = "var_2_input_object_1"
$ws.Range("C1").Value = "var_3_input_object_2"

# Update data row
$ws.Range("A2").Value = "bradleystewart"
$ws.Range("B2").Value = "Sample text"
$ws.Range("C2").Value = "ejennings"
